$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 175.411433
$ws.Range("H2").Value = 526.234299
$ws.Range("I2").Value = 0.1535106429347505
$ws.Range("J2").Value = 0.1535106429347505
$ws.Range("M2").Value = 16.790963
$ws.Range("N2").Value = 50.372889
$ws.Range("O2").Value = 0.3767143164125142
$ws.Range("P2").Value = 0.3767143164125141
$ws.Range("Q2").Value = 2945.326881279979
$ws.Range("R2").Value = 26507.94193151981
$ws.Range("S2").Value = 0.05782965691521007
$ws.Range("T2").Value = 0.05782965691521009
$ws.Range("G3").Value = 175.411433
$ws.Range("H3").Value = 526.234299
$ws.Range("I3").Value = 0.1535106429347505
$ws.Range("J3").Value = 0.1535106429347505
$ws.Range("O3").Value = 0.05559285173193915
$ws.Range("P3").Value = 0.05559285173193915
$ws.Range("Q3").Value = 434.6506450097136
$ws.Range("R3").Value = 3911.855805087423
$ws.Range("S3").Value = 0.008534094411946237
$ws.Range("T3").Value = 0.008534094411946238
$ws.Range("G4").Value = 175.411433
$ws.Range("H4").Value = 526.234299
$ws.Range("I4").Value = 0.1535106429347505
$ws.Range("J4").Value = 0.1535106429347505
$ws.Range("M4").Value = 7.558934333333333
$ws.Range("N4").Value = 22.676803
$ws.Range("O4").Value = 0.1695887710662426
$ws.Range("P4").Value = 0.1695887710662426
$ws.Range("Q4").Value = 1325.9235033629
$ws.Range("R4").Value = 11933.3115302661
$ws.Range("S4").Value = 0.02603368128089311
$ws.Range("T4").Value = 0.02603368128089312
$ws.Range("G5").Value = 175.411433
$ws.Range("H5").Value = 526.234299
$ws.Range("I5").Value = 0.1535106429347505
$ws.Range("J5").Value = 0.1535106429347505
$ws.Range("M5").Value = 17.74434966666666
$ws.Range("N5").Value = 53.23304899999999
$ws.Range("O5").Value = 0.3981040607893041
$ws.Range("P5").Value = 0.398104060789304
$ws.Range("Q5").Value = 3112.561802683072
$ws.Range("R5").Value = 28013.05622414765
$ws.Range("S5").Value = 0.06111321032670106
$ws.Range("T5").Value = 0.06111321032670107
$ws.Range("I6").Value = 0.464799214434963
$ws.Range("J6").Value = 0.4647992144349631
$ws.Range("M6").Value = 16.790963
$ws.Range("N6").Value = 50.372889
$ws.Range("O6").Value = 0.3767143164125142
$ws.Range("P6").Value = 0.3767143164125141
$ws.Range("Q6").Value = 8917.854778674853
$ws.Range("R6").Value = 80260.69300807368
$ws.Range("S6").Value = 0.1750965183349407
$ws.Range("T6").Value = 0.1750965183349407
$ws.Range("I7").Value = 0.464799214434963
$ws.Range("J7").Value = 0.4647992144349631
$ws.Range("O7").Value = 0.05559285173193915
$ws.Range("P7").Value = 0.05559285173193915
$ws.Range("S7").Value = 0.02583951381320469
$ws.Range("T7").Value = 0.02583951381320469
$ws.Range("I8").Value = 0.464799214434963
$ws.Range("J8").Value = 0.4647992144349631
$ws.Range("M8").Value = 7.558934333333333
$ws.Range("N8").Value = 22.676803
$ws.Range("O8").Value = 0.1695887710662426
$ws.Range("P8").Value = 0.1695887710662426
$ws.Range("Q8").Value = 4014.628503809226
$ws.Range("R8").Value = 36131.65653428304
$ws.Range("S8").Value = 0.07882472756858035
$ws.Range("T8").Value = 0.07882472756858036
$ws.Range("I9").Value = 0.464799214434963
$ws.Range("J9").Value = 0.4647992144349631
$ws.Range("M9").Value = 17.74434966666666
$ws.Range("N9").Value = 53.23304899999999
$ws.Range("O9").Value = 0.3981040607893041
$ws.Range("P9").Value = 0.398104060789304
$ws.Range("Q9").Value = 9424.208335719688
$ws.Range("R9").Value = 84817.87502147719
$ws.Range("S9").Value = 0.1850384547182373
$ws.Range("T9").Value = 0.1850384547182373
$ws.Range("G10").Value = 360.115397
$ws.Range("H10").Value = 1080.346191
$ws.Range("I10").Value = 0.3151536087398187
$ws.Range("J10").Value = 0.3151536087398188
$ws.Range("M10").Value = 16.790963
$ws.Range("N10").Value = 50.372889
$ws.Range("O10").Value = 0.3767143164125142
$ws.Range("P10").Value = 0.3767143164125141
$ws.Range("Q10").Value = 6046.684306757312
$ws.Range("R10").Value = 54420.15876081581
$ws.Range("S10").Value = 0.1187228762813578
$ws.Range("T10").Value = 0.1187228762813578
$ws.Range("G11").Value = 360.115397
$ws.Range("H11").Value = 1080.346191
$ws.Range("I11").Value = 0.3151536087398187
$ws.Range("J11").Value = 0.3151536087398188
$ws.Range("O11").Value = 0.05559285173193915
$ws.Range("P11").Value = 0.05559285173193915
$ws.Range("Q11").Value = 892.3271813415897
$ws.Range("R11").Value = 8030.944632074307
$ws.Range("S11").Value = 0.01752028784345831
$ws.Range("T11").Value = 0.01752028784345831
$ws.Range("G12").Value = 360.115397
$ws.Range("H12").Value = 1080.346191
$ws.Range("I12").Value = 0.3151536087398187
$ws.Range("J12").Value = 0.3151536087398188
$ws.Range("M12").Value = 7.558934333333333
$ws.Range("N12").Value = 22.676803
$ws.Range("O12").Value = 0.1695887710662426
$ws.Range("P12").Value = 0.1695887710662426
$ws.Range("Q12").Value = 2722.088638345264
$ws.Range("R12").Value = 24498.79774510737
$ws.Range("S12").Value = 0.05344651320327731
$ws.Range("T12").Value = 0.05344651320327732
$ws.Range("G13").Value = 360.115397
$ws.Range("H13").Value = 1080.346191
$ws.Range("I13").Value = 0.3151536087398187
$ws.Range("J13").Value = 0.3151536087398188
$ws.Range("M13").Value = 17.74434966666666
$ws.Range("N13").Value = 53.23304899999999
$ws.Range("O13").Value = 0.3981040607893041
$ws.Range("P13").Value = 0.398104060789304
$ws.Range("Q13").Value = 6390.013524718484
$ws.Range("R13").Value = 57510.12172246636
$ws.Range("S13").Value = 0.1254639314117253
$ws.Range("T13").Value = 0.1254639314117254
$ws.Range("G14").Value = 76.02905266666666
$ws.Range("H14").Value = 228.087158
$ws.Range("I14").Value = 0.06653653389046771
$ws.Range("J14").Value = 0.06653653389046772
$ws.Range("M14").Value = 16.790963
$ws.Range("N14").Value = 50.372889
$ws.Range("O14").Value = 0.3767143164125142
$ws.Range("P14").Value = 0.3767143164125141
$ws.Range("Q14").Value = 1276.601010251051
$ws.Range("R14").Value = 11489.40909225946
$ws.Range("S14").Value = 0.02506526488100563
$ws.Range("T14").Value = 0.02506526488100563
$ws.Range("G15").Value = 76.02905266666666
$ws.Range("H15").Value = 228.087158
$ws.Range("I15").Value = 0.06653653389046771
$ws.Range("J15").Value = 0.06653653389046772
$ws.Range("O15").Value = 0.05559285173193915
$ws.Range("P15").Value = 0.05559285173193915
$ws.Range("Q15").Value = 188.3918067133295
$ws.Range("R15").Value = 1695.526260419966
$ws.Range("S15").Value = 0.003698955663329916
$ws.Range("T15").Value = 0.003698955663329916
$ws.Range("G16").Value = 76.02905266666666
$ws.Range("H16").Value = 228.087158
$ws.Range("I16").Value = 0.06653653389046771
$ws.Range("J16").Value = 0.06653653389046772
$ws.Range("M16").Value = 7.558934333333333
$ws.Range("N16").Value = 22.676803
$ws.Range("O16").Value = 0.1695887710662426
$ws.Range("P16").Value = 0.1695887710662426
$ws.Range("Q16").Value = 574.6986165328748
$ws.Range("R16").Value = 5172.287548795874
$ws.Range("S16").Value = 0.01128384901349182
$ws.Range("T16").Value = 0.01128384901349182
$ws.Range("G17").Value = 76.02905266666666
$ws.Range("H17").Value = 228.087158
$ws.Range("I17").Value = 0.06653653389046771
$ws.Range("J17").Value = 0.06653653389046772
$ws.Range("M17").Value = 17.74434966666666
$ws.Range("N17").Value = 53.23304899999999
$ws.Range("O17").Value = 0.3981040607893041
$ws.Range("P17").Value = 0.398104060789304
$ws.Range("Q17").Value = 1349.086095342749
$ws.Range("R17").Value = 12141.77485808474
$ws.Range("S17").Value = 0.02648846433264035
$ws.Range("T17").Value = 0.02648846433264035